# "Simulated Wild Card round and logged it"
# Update Rushing and Receiving stat sheets with the results of the
# simulated Wild Card round game.

$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# A.Dalton (row 2): 3DATT, RZATT
$rushing.Range("E2").Value = 3
$rushing.Range("F2").Value = 4

# D.Montgomery (row 5): 1DATT, 2DATT, 3DATT, RZATT
$rushing.Range("C5").Value = 162
$rushing.Range("D5").Value = 83
$rushing.Range("E5").Value = 32
$rushing.Range("F5").Value = 38

# R.Nall (row 7): 1DATT, 2DATT
$rushing.Range("C7").Value = 17
$rushing.Range("D7").Value = 13

# --- Receiving sheet ---------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# D.Montgomery (row 2): Short Target, Short Comp
$receiving.Range("C2").Value = 60
$receiving.Range("D2").Value = 51

# D.Williams (row 3): Short Target, Short Comp, Deep Target, Deep Comp, RZ Target, RZ Comp
$receiving.Range("C3").Value = 15
$receiving.Range("D3").Value = 12
$receiving.Range("E3").Value = 4
$receiving.Range("F3").Value = 2
$receiving.Range("G3").Value = 3
$receiving.Range("H3").Value = 3

# D.Mooney (row 5): Short Target, Short Comp
$receiving.Range("C5").Value = 50
$receiving.Range("D5").Value = 33

# M.Goodwin (row 6): Short Target, Short Comp, Deep Target, RZ Target, RZ Comp
$receiving.Range("C6").Value = 105
$receiving.Range("D6").Value = 67
$receiving.Range("E6").Value = 35
$receiving.Range("G6").Value = 11
$receiving.Range("H6").Value = 6

# J.Grant (row 7): Short Target, Short Comp
$receiving.Range("C7").Value = 29
$receiving.Range("D7").Value = 15

# D.Byrd (row 8): Short Target, Short Comp
$receiving.Range("C8").Value = 32
$receiving.Range("D8").Value = 23

# D.Newsome (row 10): Short Target, Short Comp, Deep Target
$receiving.Range("C10").Value = 3
$receiving.Range("D10").Value = 2
$receiving.Range("E10").Value = 1

# C.Kmet (row 11): Short Target, Short Comp, Deep Target, Deep Comp
$receiving.Range("C11").Value = 80
$receiving.Range("D11").Value = 53
$receiving.Range("E11").Value = 12
$receiving.Range("F11").Value = 7

# J.Graham (row 12): Short Target, Short Comp
$receiving.Range("C12").Value = 20
$receiving.Range("D12").Value = 12
